$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new weekly record as row 91 - this shifts the existing
# rows 91..103 down to 92..104 (row 104 therefore becomes a duplicate of
# the old row 103, as in the target diff).
$ws.Rows.Item(91).Insert()

$ws.Range("A91").Value = 7
$ws.Range("B91").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C91").Value = "Ñuble"
$ws.Range("D91").Value = Get-Date -Year 2023 -Month 7 -Day 20 -Hour 0 -Minute 0 -Second 0
$ws.Range("E91").Value = 16
$ws.Range("F91").Value = 100112044
$ws.Range("G91").Value = "Perejil"
$ws.Range("H91").Value = "Sin especificar"
$ws.Range("I91").Value = "Primera"
$ws.Range("J91").Value = 100
$ws.Range("K91").Value = 1500
$ws.Range("L91").Value = 1500
$ws.Range("M91").Value = 1500
$ws.Range("N91").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O91").Value = "Región de Ñuble"
$ws.Range("P91").Value = 1500
$ws.Range("Q91").Value = 1
$ws.Range("R91").Value = "Hortaliza"
